$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update DAMSLTag (column I) and DialogAct (column J) values per re-run of SGNN dialog-act tagger
$ws.Cells.Item(2, 9).Value = 'sv'
$ws.Cells.Item(2, 10).Value = 'Statement-opinion'
$ws.Cells.Item(9, 9).Value = 'sd'
$ws.Cells.Item(9, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(13, 9).Value = 'ba'
$ws.Cells.Item(13, 10).Value = 'Appreciation'
$ws.Cells.Item(20, 9).Value = 'sv'
$ws.Cells.Item(20, 10).Value = 'Statement-opinion'
$ws.Cells.Item(30, 9).Value = 'aa'
$ws.Cells.Item(30, 10).Value = 'Agree/Accept'
$ws.Cells.Item(32, 9).Value = 'aa'
$ws.Cells.Item(32, 10).Value = 'Agree/Accept'
$ws.Cells.Item(40, 9).Value = 'aa'
$ws.Cells.Item(40, 10).Value = 'Agree/Accept'
$ws.Cells.Item(50, 9).Value = 'sd'
$ws.Cells.Item(50, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(80, 9).Value = 'aa'
$ws.Cells.Item(80, 10).Value = 'Agree/Accept'
$ws.Cells.Item(81, 9).Value = '%'
$ws.Cells.Item(81, 10).Value = 'Uninterpretable'
$ws.Cells.Item(82, 9).Value = 'aa'
$ws.Cells.Item(82, 10).Value = 'Agree/Accept'
$ws.Cells.Item(83, 9).Value = '%'
$ws.Cells.Item(83, 10).Value = 'Uninterpretable'
$ws.Cells.Item(101, 9).Value = 'sd'
$ws.Cells.Item(101, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(107, 9).Value = 'sd'
$ws.Cells.Item(107, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(113, 9).Value = 'sd'
$ws.Cells.Item(113, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(114, 9).Value = 'aa'
$ws.Cells.Item(114, 10).Value = 'Agree/Accept'
$ws.Cells.Item(124, 9).Value = 'sd'
$ws.Cells.Item(124, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(125, 9).Value = 'sd'
$ws.Cells.Item(125, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(128, 9).Value = 'sd'
$ws.Cells.Item(128, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(166, 9).Value = 'sd'
$ws.Cells.Item(166, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(174, 9).Value = '%'
$ws.Cells.Item(174, 10).Value = 'Uninterpretable'
$ws.Cells.Item(179, 9).Value = 'ba'
$ws.Cells.Item(179, 10).Value = 'Appreciation'
$ws.Cells.Item(188, 9).Value = 'aa'
$ws.Cells.Item(188, 10).Value = 'Agree/Accept'
$ws.Cells.Item(205, 9).Value = 'sv'
$ws.Cells.Item(205, 10).Value = 'Statement-opinion'
$ws.Cells.Item(209, 9).Value = 'b'
$ws.Cells.Item(209, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(211, 9).Value = 'aa'
$ws.Cells.Item(211, 10).Value = 'Agree/Accept'
